# Framer_Upload_template.xlsx update:
# Adds new "status"/manpower/registration related columns (P:Y) to the
# template header rows and refreshes column sizing / selection, matching
# the author's "Updated template file for farmer" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
# Inserted in the same order the author originally typed/pasted them in
# (status block first, then the manpower/season block, then the
# ancien/registration block) so the shared-string table grows the same
# way it did in the real edit.
$ws.Range("X1").Value = "status"
$ws.Range("X2").Value = "PENDING/APPROVED/EXCLUDED"
$ws.Range("Y1").Value = "statusComment"

$ws.Range("P1").Value = "tempManpower"
$ws.Range("Q1").Value = "permanentManpower"
$ws.Range("R1").Value = "hhMembers"
$ws.Range("S1").Value = "xsaison_last"
$ws.Range("T1").Value = "xsaison_last_but_one"
$ws.Range("U1").Value = "xsaison_last_but_two"
$ws.Range("V1").Value = "ancienCode"
$ws.Range("W1").Value = "registrationStatus"

$ws.Range("W2").Value = "ANCIEN/NOUVEAU"

# --- Column width refresh ---------------------------------------------
# Existing columns get re-measured slightly (a normal side effect of
# resaving the workbook) and the new columns P:Y get sized as well.
$ws.Columns("A").ColumnWidth = 10.833333333333332
$ws.Columns("B").ColumnWidth = 12.833333333333332
$ws.Columns("E").ColumnWidth = 12.5
$ws.Columns("F").ColumnWidth = 14.666666666666668
$ws.Columns("G").ColumnWidth = 17.666666666666664
$ws.Columns("H").ColumnWidth = 13.0
$ws.Columns("I").ColumnWidth = 15.666666666666668
$ws.Columns("J").ColumnWidth = 26.333333333333336
$ws.Columns("K").ColumnWidth = 15.666666666666668
$ws.Columns("L").ColumnWidth = 14.333333333333332
$ws.Columns("M").ColumnWidth = 25.666666666666664
$ws.Columns("N").ColumnWidth = 12.666666666666668
$ws.Columns("O:S").ColumnWidth = 15.833333333333332
$ws.Columns("T").ColumnWidth = 22.666666666666664
$ws.Columns("U:W").ColumnWidth = 22.833333333333336
$ws.Columns("X").ColumnWidth = 24.333333333333336
$ws.Columns("Y").ColumnWidth = 15.666666666666668

# --- Row height ---------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16

# --- View state ----------------------------------------------------------
# Scroll toward the newly added columns and leave the selection where the
# author left it.
$excel.Goto($ws.Range("O1"), $true)
$ws.Range("T8").Select()
